# Auto-generated cell updates reproducing the sheet refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 595.9091
$ws.Range("I2").Value = 173.22223
$ws.Range("K2").Value = 173.22223
$ws.Range("M2").Value = -60.22223
$ws.Range("H43").Value = 5162.375
$ws.Range("I43").Value = 5954.222
$ws.Range("J43").Value = 4144.2856
$ws.Range("K43").Value = 5954.222
$ws.Range("L43").Value = 4144.2856
$ws.Range("M43").Value = -5885.222
$ws.Range("N43").Value = -4282.2856
$ws.Range("H62").Value = 3266.3572
$ws.Range("I62").Value = 2933.3684
$ws.Range("K62").Value = 2933.3684
$ws.Range("M62").Value = -2309.3684
$ws.Range("H65").Value = 3266.3572
$ws.Range("I65").Value = 2933.3684
$ws.Range("K65").Value = 14666.842
$ws.Range("M65").Value = -11546.842
$ws.Range("H100").Value = 3334
$ws.Range("J100").Value = 1323.5
$ws.Range("L100").Value = 1323.5
$ws.Range("N100").Value = -2405.5
$ws.Range("H107").Value = 3319.9285
$ws.Range("I107").Value = 1644
$ws.Range("K107").Value = 1644
$ws.Range("M107").Value = 276
$ws.Range("H112").Value = 4250.4
$ws.Range("J112").Value = 3927.5
$ws.Range("L112").Value = 11782.5
$ws.Range("N112").Value = -13998.5
$ws.Range("H113").Value = 6591.467
$ws.Range("I113").Value = 6231.222
$ws.Range("K113").Value = 6231.222
$ws.Range("M113").Value = -2977.222
$ws.Range("H116").Value = 17327.65
$ws.Range("I116").Value = 14159.5625
$ws.Range("K116").Value = 14159.5625
$ws.Range("M116").Value = -10717.5625
$ws.Range("H127").Value = 2140
$ws.Range("J127").Value = 2453.875
$ws.Range("L127").Value = 7361.625
$ws.Range("N127").Value = -17281.625
$ws.Range("H131").Value = 1626192.6
$ws.Range("I131").Value = 1922.7
$ws.Range("J131").Value = 5686867.5
$ws.Range("K131").Value = 5768.1
$ws.Range("L131").Value = 17060602.5
$ws.Range("M131").Value = -728.1000000000004
$ws.Range("N131").Value = -17070682.5
$ws.Range("H132").Value = 1497.8611
$ws.Range("I132").Value = 1338.1562
$ws.Range("K132").Value = 4014.4686
$ws.Range("M132").Value = -1484.4686
$ws.Range("H137").Value = 2649.5386
$ws.Range("I137").Value = 1874.75
$ws.Range("J137").Value = 3889.2
$ws.Range("K137").Value = 5624.25
$ws.Range("L137").Value = 11667.6
$ws.Range("M137").Value = -3074.25
$ws.Range("N137").Value = -16767.6
$ws.Range("H141").Value = 3539.7837
$ws.Range("I141").Value = 2713.1177
$ws.Range("K141").Value = 8139.353099999999
$ws.Range("M141").Value = -2959.353099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1186.4445
$ws.Range("I2").Value = 1147.25
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 1147.25
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -1034.25
$ws.Range("N2").Value = -1726
$ws.Range("H32").Value = 13633.593
$ws.Range("I32").Value = 11569.913
$ws.Range("K32").Value = 11569.913
$ws.Range("M32").Value = -11282.913
$ws.Range("H45").Value = 2927.5454
$ws.Range("I45").Value = 1562.3
$ws.Range("J45").Value = 4065.25
$ws.Range("K45").Value = 1562.3
$ws.Range("L45").Value = 4065.25
$ws.Range("M45").Value = -1185.3
$ws.Range("N45").Value = -4819.25
$ws.Range("H61").Value = 11059382
$ws.Range("I61").Value = 16672964
$ws.Range("J61").Value = 1436099.6
$ws.Range("K61").Value = 16672964
$ws.Range("L61").Value = 1436099.6
$ws.Range("M61").Value = -16672752
$ws.Range("N61").Value = -1436523.6
$ws.Range("H74").Value = 3395.8667
$ws.Range("I74").Value = 3281.3572
$ws.Range("J74").Value = 4999
$ws.Range("K74").Value = 3281.3572
$ws.Range("L74").Value = 4999
$ws.Range("M74").Value = -2407.3572
$ws.Range("N74").Value = -6747
$ws.Range("H77").Value = 3395.8667
$ws.Range("I77").Value = 3281.3572
$ws.Range("J77").Value = 4999
$ws.Range("K77").Value = 16406.786
$ws.Range("L77").Value = 24995
$ws.Range("M77").Value = -12038.786
$ws.Range("N77").Value = -33731
$ws.Range("H110").Value = 9725.6
$ws.Range("J110").Value = 4606.7144
$ws.Range("L110").Value = 4606.7144
$ws.Range("N110").Value = -8696.714400000001
$ws.Range("H116").Value = 1186.4445
$ws.Range("I116").Value = 1147.25
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 1147.25
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 1146.75
$ws.Range("N116").Value = -6088
$ws.Range("H122").Value = 5299.1816
$ws.Range("I122").Value = 5143.4443
$ws.Range("K122").Value = 15430.3329
$ws.Range("M122").Value = -12980.3329
$ws.Range("H132").Value = 3705988.8
$ws.Range("I132").Value = 2581
$ws.Range("J132").Value = 16667916
$ws.Range("K132").Value = 7743
$ws.Range("L132").Value = 50003748
$ws.Range("M132").Value = -5213
$ws.Range("N132").Value = -50008808
$ws.Range("H133").Value = 93001
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 11059382
$ws.Range("I136").Value = 16672964
$ws.Range("J136").Value = 1436099.6
$ws.Range("K136").Value = 50018892
$ws.Range("L136").Value = 4308298.800000001
$ws.Range("M136").Value = -50016342
$ws.Range("N136").Value = -4313398.800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1186.4445
$ws.Range("I3").Value = 1147.25
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 1147.25
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -1033.25
$ws.Range("N3").Value = -1728
$ws.Range("H20").Value = 3076.1482
$ws.Range("I20").Value = 2630.1765
$ws.Range("J20").Value = 3834.3
$ws.Range("K20").Value = 2630.1765
$ws.Range("L20").Value = 3834.3
$ws.Range("M20").Value = -2383.1765
$ws.Range("N20").Value = -4328.3
$ws.Range("H81").Value = 32556
$ws.Range("J81").Value = 32556
$ws.Range("L81").Value = 32556
$ws.Range("N81").Value = -34678
$ws.Range("H84").Value = 32556
$ws.Range("J84").Value = 32556
$ws.Range("L84").Value = 97668
$ws.Range("N84").Value = -108276
$ws.Range("H86").Value = 5121.7334
$ws.Range("I86").Value = 3649.5
$ws.Range("J86").Value = 6804.2856
$ws.Range("K86").Value = 3649.5
$ws.Range("L86").Value = 6804.2856
$ws.Range("M86").Value = -2526.5
$ws.Range("N86").Value = -9050.285599999999
$ws.Range("H89").Value = 5121.7334
$ws.Range("I89").Value = 3649.5
$ws.Range("J89").Value = 6804.2856
$ws.Range("K89").Value = 18247.5
$ws.Range("L89").Value = 34021.428
$ws.Range("M89").Value = -12631.5
$ws.Range("N89").Value = -45253.428
$ws.Range("H99").Value = 2082.6365
$ws.Range("I99").Value = 2221
$ws.Range("K99").Value = 2221
$ws.Range("M99").Value = -723
$ws.Range("H102").Value = 32888.25
$ws.Range("I102").Value = 10518
$ws.Range("K102").Value = 10518
$ws.Range("M102").Value = -7273
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H134").Value = 2858672.8
$ws.Range("I134").Value = 1575
$ws.Range("K134").Value = 4725
$ws.Range("M134").Value = -2190
$ws.Range("H135").Value = 82090.37
$ws.Range("J135").Value = 82090.37
$ws.Range("L135").Value = 82090.37
$ws.Range("N135").Value = -92230.37

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1365.375
$ws.Range("I22").Value = 570.5
$ws.Range("K22").Value = 570.5
$ws.Range("M22").Value = -220.5
$ws.Range("H31").Value = 40005976
$ws.Range("I31").Value = 111116630
$ws.Range("K31").Value = 111116630
$ws.Range("M31").Value = -111116335
$ws.Range("H34").Value = 40005976
$ws.Range("I34").Value = 111116630
$ws.Range("K34").Value = 111116630
$ws.Range("M34").Value = -111116428
$ws.Range("H52").Value = 84999
$ws.Range("J52").Value = 84999
$ws.Range("L52").Value = 84999
$ws.Range("N52").Value = -85587
$ws.Range("H86").Value = 6031.769
$ws.Range("J86").Value = 6415.8335
$ws.Range("L86").Value = 6415.8335
$ws.Range("N86").Value = -8661.833500000001
$ws.Range("H89").Value = 6031.769
$ws.Range("J89").Value = 6415.8335
$ws.Range("L89").Value = 32079.1675
$ws.Range("N89").Value = -43311.1675
$ws.Range("H102").Value = 73309.75
$ws.Range("J102").Value = 92746.336
$ws.Range("L102").Value = 92746.336
$ws.Range("N102").Value = -97614.336
$ws.Range("H103").Value = 42380
$ws.Range("I103").Value = 15594
$ws.Range("K103").Value = 15594
$ws.Range("M103").Value = -14422
$ws.Range("H107").Value = 857.70966
$ws.Range("I107").Value = 448.42105
$ws.Range("K107").Value = 448.42105
$ws.Range("M107").Value = 1471.57895
$ws.Range("H116").Value = 99999
$ws.Range("J116").Value = 99999
$ws.Range("L116").Value = 99999
$ws.Range("N116").Value = -109177
$ws.Range("H118").Value = 99998
$ws.Range("J118").Value = 99998
$ws.Range("L118").Value = 99998
$ws.Range("N118").Value = -103312
$ws.Range("H119").Value = 63166
$ws.Range("J119").Value = 63166
$ws.Range("L119").Value = 63166
$ws.Range("N119").Value = -72842
$ws.Range("H122").Value = 2093.4
$ws.Range("I122").Value = 2093.4
$ws.Range("K122").Value = 6280.200000000001
$ws.Range("M122").Value = -3830.200000000001
$ws.Range("H129").Value = 74199.5
$ws.Range("J129").Value = 74199.5
$ws.Range("L129").Value = 74199.5
$ws.Range("N129").Value = -84199.5
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 87499
$ws.Range("J131").Value = 87499
$ws.Range("L131").Value = 87499
$ws.Range("N131").Value = -97579
$ws.Range("H132").Value = 1669.25
$ws.Range("I132").Value = 1693.9474
$ws.Range("K132").Value = 5081.8422
$ws.Range("M132").Value = -2551.8422
$ws.Range("H141").Value = 514806.75
$ws.Range("I141").Value = 389998
$ws.Range("J141").Value = 532636.5600000001
$ws.Range("K141").Value = 389998
$ws.Range("L141").Value = 532636.5600000001
$ws.Range("M141").Value = -384818
$ws.Range("N141").Value = -542996.5600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 6479
$ws.Range("I9").Value = 1250
$ws.Range("J9").Value = 22166
$ws.Range("K9").Value = 3750
$ws.Range("L9").Value = 66498
$ws.Range("M9").Value = -3526
$ws.Range("N9").Value = -66946
$ws.Range("H69").Value = 16714
$ws.Range("I69").Value = 3094.75
$ws.Range("J69").Value = 23523.625
$ws.Range("K69").Value = 9284.25
$ws.Range("L69").Value = 70570.875
$ws.Range("M69").Value = -8473.25
$ws.Range("N69").Value = -72192.875
$ws.Range("H72").Value = 16714
$ws.Range("I72").Value = 3094.75
$ws.Range("J72").Value = 23523.625
$ws.Range("K72").Value = 27852.75
$ws.Range("L72").Value = 211712.625
$ws.Range("M72").Value = -23796.75
$ws.Range("N72").Value = -219824.625
$ws.Range("H82").Value = 21124
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 21124
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 63372
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -64184
$ws.Range("H85").Value = 21124
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 21124
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 63372
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -66180
$ws.Range("H92").Value = 69.75
$ws.Range("J92").Value = 69.75
$ws.Range("L92").Value = 209.25
$ws.Range("N92").Value = -2705.25
$ws.Range("H97").Value = 751.5
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 751.5
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2254.5
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -3246.5
$ws.Range("H125").Value = 22932.2
$ws.Range("J125").Value = 33332
$ws.Range("L125").Value = 99996
$ws.Range("N125").Value = -109836
$ws.Range("H131").Value = 4372.9473
$ws.Range("J131").Value = 5252.5
$ws.Range("L131").Value = 15757.5
$ws.Range("N131").Value = -25837.5
$ws.Range("H138").Value = 16542.846
$ws.Range("I138").Value = 14064.8
$ws.Range("K138").Value = 42194.39999999999
$ws.Range("M138").Value = -37054.39999999999
$ws.Range("H139").Value = 6015
$ws.Range("I139").Value = 1532.1111
$ws.Range("K139").Value = 4596.3333
$ws.Range("M139").Value = 543.6666999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 71.61539
$ws.Range("I2").Value = 72.59999999999999
$ws.Range("K2").Value = 72.59999999999999
$ws.Range("M2").Value = 40.40000000000001
$ws.Range("H39").Value = 39999.5
$ws.Range("J39").Value = 39999.5
$ws.Range("L39").Value = 39999.5
$ws.Range("N39").Value = -41063.5
$ws.Range("H102").Value = 1696.92
$ws.Range("I102").Value = 1743.3914
$ws.Range("K102").Value = 1743.3914
$ws.Range("M102").Value = -121.3914
$ws.Range("H107").Value = 1175.5
$ws.Range("I107").Value = 678.4286
$ws.Range("K107").Value = 678.4286
$ws.Range("M107").Value = 1241.5714
$ws.Range("H113").Value = 1326468.8
$ws.Range("J113").Value = 3090041.2
$ws.Range("L113").Value = 3090041.2
$ws.Range("N113").Value = -3094381.2
$ws.Range("H122").Value = 6273
$ws.Range("I122").Value = 8924.5
$ws.Range("J122").Value = 3621.5
$ws.Range("K122").Value = 26773.5
$ws.Range("L122").Value = 10864.5
$ws.Range("M122").Value = -24323.5
$ws.Range("N122").Value = -15764.5
$ws.Range("H132").Value = 3033712.8
$ws.Range("I132").Value = 3634.1853
$ws.Range("J132").Value = 16669066
$ws.Range("K132").Value = 10902.5559
$ws.Range("L132").Value = 50007198
$ws.Range("M132").Value = -8372.555899999999
$ws.Range("N132").Value = -50012258

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7977
$ws.Range("I7").Value = 7977
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 7977
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -7865
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 4779
$ws.Range("I22").Value = 4724
$ws.Range("J22").Value = 4999
$ws.Range("K22").Value = 4724
$ws.Range("L22").Value = 4999
$ws.Range("M22").Value = -4429
$ws.Range("N22").Value = -5589
$ws.Range("H27").Value = 4779
$ws.Range("I27").Value = 4724
$ws.Range("J27").Value = 4999
$ws.Range("K27").Value = 4724
$ws.Range("L27").Value = 4999
$ws.Range("M27").Value = -4617
$ws.Range("N27").Value = -5213
$ws.Range("H46").Value = 1884.3636
$ws.Range("I46").Value = 1887.8
$ws.Range("J46").Value = 1881.5
$ws.Range("K46").Value = 1887.8
$ws.Range("L46").Value = 1881.5
$ws.Range("M46").Value = -1699.8
$ws.Range("N46").Value = -2257.5
$ws.Range("H47").Value = 590
$ws.Range("I47").Value = 590
$ws.Range("K47").Value = 590
$ws.Range("M47").Value = -100
$ws.Range("H51").Value = 99999
$ws.Range("J51").Value = 99999
$ws.Range("L51").Value = 99999
$ws.Range("N51").Value = -100955
$ws.Range("H52").Value = 590
$ws.Range("I52").Value = 590
$ws.Range("K52").Value = 590
$ws.Range("M52").Value = -357
$ws.Range("H55").Value = 1175.5834
$ws.Range("I55").Value = 702.3333
$ws.Range("K55").Value = 702.3333
$ws.Range("M55").Value = -529.3333
$ws.Range("H57").Value = 32207.46
$ws.Range("I57").Value = 26558.166
$ws.Range("K57").Value = 26558.166
$ws.Range("M57").Value = -25992.166
$ws.Range("H109").Value = 76546
$ws.Range("J109").Value = 89994.5
$ws.Range("L109").Value = 89994.5
$ws.Range("N109").Value = -92768.5
$ws.Range("H112").Value = 109718.664
$ws.Range("J112").Value = 109718.664
$ws.Range("L112").Value = 109718.664
$ws.Range("N112").Value = -112672.664
$ws.Range("H126").Value = 7977
$ws.Range("I126").Value = 7977
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 23931
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -21461
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3975.743
$ws.Range("I132").Value = 2636
$ws.Range("K132").Value = 7908
$ws.Range("M132").Value = -5378
$ws.Range("H134").Value = 92333
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 92333
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 92333
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -102473

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 61124.75
$ws.Range("J46").Value = 61124.75
$ws.Range("L46").Value = 61124.75
$ws.Range("N46").Value = -61586.75
$ws.Range("H126").Value = 2739.3076
$ws.Range("I126").Value = 2782.818
$ws.Range("K126").Value = 8348.454000000002
$ws.Range("M126").Value = -5878.454000000002
$ws.Range("H132").Value = 402224
$ws.Range("I132").Value = 2089.4736
$ws.Range("K132").Value = 6268.4208
$ws.Range("M132").Value = -3738.4208
$ws.Range("H133").Value = 68864.336
$ws.Range("J133").Value = 68864.336
$ws.Range("L133").Value = 68864.336
$ws.Range("N133").Value = -78984.336
$ws.Range("H134").Value = 61124.75
$ws.Range("J134").Value = 61124.75
$ws.Range("L134").Value = 183374.25
$ws.Range("N134").Value = -188444.25

